$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1994.2858
$ws.Range("I9").Value = 1999
$ws.Range("J9").Value = 1993.5
$ws.Range("K9").Value = 1999
$ws.Range("L9").Value = 1993.5
$ws.Range("M9").Value = -1830
$ws.Range("N9").Value = -2331.5
$ws.Range("H80").Value = 1453.5714
$ws.Range("I80").Value = 625
$ws.Range("J80").Value = 1785
$ws.Range("K80").Value = 1875
$ws.Range("L80").Value = 5355
$ws.Range("M80").Value = -877
$ws.Range("N80").Value = -7351
$ws.Range("H83").Value = 1453.5714
$ws.Range("I83").Value = 625
$ws.Range("J83").Value = 1785
$ws.Range("K83").Value = 5625
$ws.Range("L83").Value = 16065
$ws.Range("M83").Value = -633
$ws.Range("N83").Value = -26049
$ws.Range("H116").Value = 5896.636
$ws.Range("J116").Value = 6171.143
$ws.Range("L116").Value = 6171.143
$ws.Range("N116").Value = -13055.143
$ws.Range("H138").Value = 3226.8235
$ws.Range("J138").Value = 3859.8
$ws.Range("L138").Value = 11579.4
$ws.Range("N138").Value = -21859.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5698.5557
$ws.Range("I32").Value = 3811.9092
$ws.Range("K32").Value = 3811.9092
$ws.Range("M32").Value = -3524.9092
$ws.Range("H63").Value = 4397.3335
$ws.Range("I63").Value = 2696.5
$ws.Range("J63").Value = 5247.75
$ws.Range("K63").Value = 2696.5
$ws.Range("L63").Value = 5247.75
$ws.Range("M63").Value = -2010.5
$ws.Range("N63").Value = -6619.75
$ws.Range("H66").Value = 4397.3335
$ws.Range("I66").Value = 2696.5
$ws.Range("J66").Value = 5247.75
$ws.Range("K66").Value = 13482.5
$ws.Range("L66").Value = 26238.75
$ws.Range("M66").Value = -10050.5
$ws.Range("N66").Value = -33102.75
$ws.Range("H132").Value = 1006
$ws.Range("I132").Value = 1006
$ws.Range("K132").Value = 3018
$ws.Range("M132").Value = -488

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3581.7
$ws.Range("J134").Value = 1997.5
$ws.Range("L134").Value = 5992.5
$ws.Range("N134").Value = -11062.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2200
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1576
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2200
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -7880
$ws.Range("N65").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 1917.7778
$ws.Range("I132").Value = 1947.2
$ws.Range("K132").Value = 5841.6
$ws.Range("M132").Value = -3311.6
$ws.Range("H134").Value = 1309.1428
$ws.Range("I134").Value = 1309.1428
$ws.Range("K134").Value = 3927.4284
$ws.Range("M134").Value = -1392.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1883.125
$ws.Range("I22").Value = 2149
$ws.Range("J22").Value = 1440
$ws.Range("K22").Value = 6447
$ws.Range("L22").Value = 4320
$ws.Range("M22").Value = -6278
$ws.Range("N22").Value = -4658
$ws.Range("H26").Value = 22.6
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = 78
$ws.Range("L26").Value = 27
$ws.Range("M26").Value = 210
$ws.Range("N26").Value = -603
$ws.Range("H27").Value = 1883.125
$ws.Range("I27").Value = 2149
$ws.Range("J27").Value = 1440
$ws.Range("K27").Value = 6447
$ws.Range("L27").Value = 4320
$ws.Range("M27").Value = -6345
$ws.Range("N27").Value = -4524
$ws.Range("H40").Value = 100.75
$ws.Range("I40").Value = 103.27273
$ws.Range("J40").Value = 95.2
$ws.Range("K40").Value = 413.09092
$ws.Range("L40").Value = 380.8
$ws.Range("M40").Value = -344.09092
$ws.Range("N40").Value = -518.8
$ws.Range("H113").Value = 294.5
$ws.Range("J113").Value = 330.25
$ws.Range("L113").Value = 990.75
$ws.Range("N113").Value = -5330.75
$ws.Range("H129").Value = 1790.8182
$ws.Range("I129").Value = 1166.1666
$ws.Range("J129").Value = 2540.4
$ws.Range("K129").Value = 3498.4998
$ws.Range("L129").Value = 7621.200000000001
$ws.Range("M129").Value = 1501.5002
$ws.Range("N129").Value = -17621.2
$ws.Range("H131").Value = 1078.5
$ws.Range("J131").Value = 1076.8462
$ws.Range("L131").Value = 3230.5386
$ws.Range("N131").Value = -13310.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 25500
$ws.Range("J47").Value = 25500
$ws.Range("L47").Value = 25500
$ws.Range("N47").Value = -26636
$ws.Range("H126").Value = 90912890
$ws.Range("J126").Value = 4137.375
$ws.Range("L126").Value = 12412.125
$ws.Range("N126").Value = -17352.125
$ws.Range("H132").Value = 1560.9166
$ws.Range("I132").Value = 1568.4546
$ws.Range("K132").Value = 4705.3638
$ws.Range("M132").Value = -2175.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3629.2856
$ws.Range("I40").Value = 3440.84
$ws.Range("K40").Value = 3440.84
$ws.Range("M40").Value = -3304.84
$ws.Range("H46").Value = 1827.2273
$ws.Range("I46").Value = 1346.6
$ws.Range("J46").Value = 2857.1428
$ws.Range("K46").Value = 1346.6
$ws.Range("L46").Value = 2857.1428
$ws.Range("M46").Value = -1158.6
$ws.Range("N46").Value = -3233.1428
$ws.Range("H68").Value = 4003
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 4003
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 5830.5557
$ws.Range("I122").Value = 5622.375
$ws.Range("J122").Value = 7496
$ws.Range("K122").Value = 16867.125
$ws.Range("L122").Value = 22488
$ws.Range("M122").Value = -14417.125
$ws.Range("N122").Value = -27388
$ws.Range("H132").Value = 5804.9165
$ws.Range("I132").Value = 6462.4443
$ws.Range("J132").Value = 3832.3333
$ws.Range("K132").Value = 19387.3329
$ws.Range("L132").Value = 11496.9999
$ws.Range("M132").Value = -16857.3329
$ws.Range("N132").Value = -16556.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 288.33334
$ws.Range("I107").Value = 308.375
$ws.Range("K107").Value = 925.125
$ws.Range("M107").Value = 994.875
$ws.Range("H113").Value = 392
$ws.Range("I113").Value = 318.33334
$ws.Range("K113").Value = 955.0000200000001
$ws.Range("M113").Value = 1214.99998
$ws.Range("H122").Value = 398.33334
$ws.Range("I122").Value = 398.33334
$ws.Range("K122").Value = 1195.00002
$ws.Range("M122").Value = 1254.99998
$ws.Range("H132").Value = 802.6957
$ws.Range("I132").Value = 611.9091
$ws.Range("K132").Value = 1835.7273
$ws.Range("M132").Value = 694.2727
$ws.Range("H136").Value = 1870.8572
$ws.Range("I136").Value = 1399.5385
$ws.Range("K136").Value = 4198.6155
$ws.Range("M136").Value = -1648.6155
